# SpellSuccessPercentage.xlsx - small oracle improvements;
# fountain appears in 1 in 20 rooms.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: Spell level (C2) 17 -> 18 ---
$ws.Range("C2").Value = 18

# --- Row 3: Penalty base/unit values ---
$ws.Range("C3").Value = 25
$ws.Range("H3").Value = 0

# --- Row 5 ---
$ws.Range("C5").Value = 21
$ws.Range("H5").Value = 0

# --- Row 6 ---
$ws.Range("C6").Value = 3

# --- Row 7 ---
$ws.Range("H7").Value = 0

# I2:I7 all share the same relative formula pattern (Skill * modifier),
# which is what makes Excel collapse them into one shared formula group.
$ws.Range("I2:I7").FormulaR1C1 = "=RC[-2]*RC[-1]"

# --- Row 11: new "Spell cost" line (F11 label, G11 = C3*3) ---
$ws.Range("F11").Value = "Spell cost"
$ws.Range("G11").Formula = "=C3*3"

# Move/restore the active selection to C5, matching the saved view state.
$ws.Range("C5").Select()

$wb.Save()
